# Update the "取得日時" (retrieved-at) timestamp in column A for all data rows
# (rows 2-10) on the "ランサーズ" sheet from 2026-01-03 01:22:01 to
# 2026-01-03 01:53:49, reflecting a re-run of the scraper at that time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-03 01:53:49"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
